# Weekly refresh of the Jengibre (Vega Modelo de Temuco) price series:
# a new week's record is inserted at row 51 (pushing the existing
# historical rows down by one), dated 2021-10-26 (serial 44495) with the
# same volume/price figures as the record that used to sit in that slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 51; Excel shifts rows 51:116
# down to 52:117 and copies row 50's formatting (incl. the date style)
# onto the new row.
$ws.Rows.Item(51).Insert()

$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 44495
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 100114007
$ws.Range("G51").Value = "Jengibre"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 30
$ws.Range("K51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = 25000
$ws.Range("N51").Value = "$/caja 13 kilos"
$ws.Range("O51").Value = "Perú"
$ws.Range("P51").Value = 1923
$ws.Range("Q51").Value = 13
$ws.Range("R51").Value = "Hortaliza"
